$wb = $excel.ActiveWorkbook

# The F column ("想去人数") values changed on the rows 2-5 for both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same events).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 898
    $ws.Range("F3").Value = 4551
    $ws.Range("F4").Value = 129
    $ws.Range("F5").Value = 794
}
